# Weekly update: insert a new "Fruta / hortaliza, semanal" record as the
# first data row for this subset (new row 8), pushing the existing rows
# 8-23 down to 9-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 8:23 down to 9:24, leaving a fresh blank row 8.
$ws.Rows("8:8").Insert()

# Populate the new row 8 with this week's observation (same market /
# product metadata as the rest of the block, new date + price figures).
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C8").Value = 'Arica y Parinacota'
$ws.Range("D8").Value = 44525
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 'Fruta'
$ws.Range("G8").Value = 100104
$ws.Range("H8").Value = 'Frutos de pepita'
$ws.Range("I8").Value = 100104005
$ws.Range("J8").Value = 'Pera'
$ws.Range("K8").Value = 'Packham''s Triumph'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1083
$ws.Range("T8").Value = 18
